$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing comment on B5 ---
$c5 = $ws.Range("B5").Comment
$c5.Text("admin:" + "`r`n" + "Belum ada di mapping, tapi sudah ada di script")

# --- Add new comments B6:B10 (same note text as updated B5) ---
$noteScript = "admin:" + "`r`n" + "Belum ada di mapping, tapi sudah ada di script"
$ws.Range("B6").AddComment($noteScript)
$ws.Range("B7").AddComment($noteScript)
$ws.Range("B8").AddComment($noteScript)
$ws.Range("B9").AddComment($noteScript)
$ws.Range("B10").AddComment($noteScript)

# --- Add new comment B13 ---
$ws.Range("B13").AddComment("admin:" + "`r`n" + "Belum ada di Mapping")

# --- Add new comment B14 ---
$ws.Range("B14").AddComment("admin:" + "`r`n" + "Belum ada di mapping tapi ada di script")

# --- Add new comments B16:B18 ---
$noteAda = "admin:" + "`r`n" + "Belum ada di mapping, tapi ada di script"
$ws.Range("B16").AddComment($noteAda)
$ws.Range("B17").AddComment($noteAda)
$ws.Range("B18").AddComment($noteAda)

# --- New rows of data (order matters for shared-string allocation) ---
$ws.Range("B14").Value = "ste_cswnrecwo"
$ws.Range("A14").Value = "matrectrans"
$ws.Range("A15").Value = "Line 665 dan 666 di sheet 0003 purchase order sepertinya salah kode package dan nama table, harusnya matrectrans bukan matusetrans"
$ws.Range("B16").Value = "ste_cswnctrycd"
$ws.Range("B17").Value = "ste_cswnctryname"
$ws.Range("B18").Value = "ste_cswndpt"
$ws.Range("A16").Value = "CURRENCY"
$ws.Range("A17").Value = "CURRENCY"
$ws.Range("A18").Value = "CURRENCY"

# --- Update selection to match the saved cursor position ---
$ws.Range("B7").Select()
